$d = $word.ActiveDocument

function Merge-IdRun([string]$idValue) {
    $searchText = "<id>" + $idValue + "</id>"

    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return
    }

    $startPos = $rng.Start
    $endPos = $rng.End

    $escaped = $idValue -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
    $combinedText = "&lt;id&gt;" + $escaped + "&lt;/id&gt;"

    # Work out whether the matched text is the entirety of its paragraph's
    # content (save for the trailing paragraph mark). When it is, we can
    # safely delete + InsertXML a fully-specified <w:p>/<w:r> fragment and
    # get an exact, attribute-faithful merged run. When the match sits in
    # the middle of a longer paragraph, InsertXML-ing a <w:p> fragment
    # there corrupts the paragraph (splits it / inserts stray breaks), so
    # we fall back to a plain Range.Text assignment, which Word performs
    # as an in-place run merge (keeping the first run's character
    # formatting) without disturbing surrounding content.
    $paraRng = $d.Range($startPos, $endPos)
    $paraRng.Expand(4) | Out-Null   # wdParagraph = 4
    $isFullParagraph = ($startPos -eq $paraRng.Start) -and ($endPos -eq ($paraRng.End - 1))

    if ($isFullParagraph) {
        $delRng = $d.Range($startPos, $endPos)
        $delRng.Delete()

        $insRng = $d.Range($startPos, $startPos)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body><w:p><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' +
            '<w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
            '<w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr>' +
            '<w:t xml:space="preserve">' + $combinedText + '</w:t>' +
            '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $insRng.InsertXML($xml)
    } else {
        # Mid-paragraph match: a plain Range.Text assignment correctly
        # merges the covered runs into one run using the first run's
        # character formatting, without disturbing the rest of the
        # paragraph - as long as Word actually perceives a text change.
        # Because the replacement text is identical to what Find just
        # matched, assigning it directly would be a no-op (the runs
        # would stay split). Force a genuine change by going through a
        # placeholder value first, then rewriting that placeholder with
        # the real (identical-looking) text, which merges the runs.
        $placeholder = "ZZZ_MERGE_PLACEHOLDER_ZZZ"

        $rngA = $d.Range($startPos, $endPos)
        $rngA.Text = $placeholder

        $rngB = $d.Range($startPos, $startPos + $placeholder.Length)
        $rngB.Text = "<id>" + $idValue + "</id>"
    }
}

Merge-IdRun "p109r_5"
Merge-IdRun "p109v_1"
Merge-IdRun "p109v_2"
